# Insert a new weekly price record for Mango at Vega Central Mapocho de
# Santiago. The new observation is inserted as row 677, pushing the
# existing rows 677-761 down to 678-762 (dimension grows from
# A1:T761 to A1:T762).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data down by inserting a fresh row at 677.
$ws.Rows.Item(677).Insert()

# Populate the newly inserted row with the latest observation.
$ws.Cells.Item(677, 1).Value = 9
$ws.Cells.Item(677, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(677, 3).Value = "Metropolitana"
$ws.Cells.Item(677, 4).Value = 45212
$ws.Cells.Item(677, 5).Value = 13
$ws.Cells.Item(677, 6).Value = "Fruta"
$ws.Cells.Item(677, 7).Value = 100108
$ws.Cells.Item(677, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(677, 9).Value = 100108002
$ws.Cells.Item(677, 10).Value = "Mango"
$ws.Cells.Item(677, 11).Value = "Sin especificar"
$ws.Cells.Item(677, 12).Value = "Primera"
$ws.Cells.Item(677, 13).Value = 700
$ws.Cells.Item(677, 14).Value = 8000
$ws.Cells.Item(677, 15).Value = 8500
$ws.Cells.Item(677, 16).Value = 8179
$ws.Cells.Item(677, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(677, 18).Value = "Brasil"
$ws.Cells.Item(677, 19).Value = 2045
$ws.Cells.Item(677, 20).Value = 4
